$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -6.862099999999997
$ws.Range("D7").Value = -8.672499999999999
$ws.Range("E7").Value = 16.06150000000001
$ws.Range("E15").Value = 16.3213
$ws.Range("D16").Value = -8.697100000000002
$ws.Range("E21").Value = 17.26739999999999
$ws.Range("E22").Value = 17.109
$ws.Range("E23").Value = 16.09399999999999
$ws.Range("D28").Value = -8.469299999999997
$ws.Range("D29").Value = -7.131899999999998
$ws.Range("D32").Value = -9.201599999999992
$ws.Range("E34").Value = 17.5881
$ws.Range("D40").Value = -7.919999999999992
$ws.Range("E43").Value = 17.50580000000001
$ws.Range("E45").Value = 16.3852
$ws.Range("E50").Value = 16.41359999999999
$ws.Range("E51").Value = 17.35450000000002
$ws.Range("D52").Value = -7.183599999999995
$ws.Range("D57").Value = -8.439399999999996
$ws.Range("D66").Value = -6.926599999999999
$ws.Range("E66").Value = 17.21280000000002
$ws.Range("E67").Value = 17.04370000000002
$ws.Range("E79").Value = 18.41420000000002
$ws.Range("E84").Value = 16.4663
$ws.Range("E92").Value = 18.58810000000001
$ws.Range("E97").Value = 16.5053
$ws.Range("D100").Value = -8.633900000000001
